$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "Förändrad" date column (C) for rows 2-6: 45207 -> 45208
$ws.Range("C2:C6").Value = 45208

# Update hyperlink formulas in row 2 (columns S-Y):
# replace "Logging_MONSTERAS" with "Logging_0861" in the URL, keep display text
$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0861/artfynd/A 32298-2023.xlsx", "A 32298-2023")'
$ws.Range("T2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0861/kartor/A 32298-2023.png", "A 32298-2023")'
$ws.Range("U2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0861/knärot/A 32298-2023.png", "A 32298-2023")'
$ws.Range("V2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0861/klagomål/A 32298-2023.docx", "A 32298-2023")'
$ws.Range("W2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0861/klagomålsmail/A 32298-2023.docx", "A 32298-2023")'
$ws.Range("X2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0861/tillsyn/A 32298-2023.docx", "A 32298-2023")'
$ws.Range("Y2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_0861/tillsynsmail/A 32298-2023.docx", "A 32298-2023")'
